# Workbook "stickers": add a new "oтвет" (File_id) column header in C1,
# and add a new data row (Пока / its sticker File_id) in row 3, matching
# the formatting already used for the "Привет" sticker row (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C (row 1)
$ws.Range("C1").Value = "oтвет"

# New data row 3: key phrase "Пока" with its sticker File_id
$ws.Range("A3").Value = "Пока"
$ws.Range("B3").Value = "CAACAgIAAxkBAAN_YiiMAbbk7OsUMHmvH_tWKG24c5cAAjQBAAJSiZEjE83Xb_UcB1gjBA"

# Match the formatting already used for the File_id cell above it (B2:
# Segoe UI 12, black) by copying its format onto the new cell.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match row 2's (taller) row height for the new row
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight

# Leave the active selection on the newly entered sticker id cell
$ws.Range("B3").Select()
